# The footer-style signature block on the cover page contains the email
# address "walter@walterkwaninsurance.com" (the text is actually split over
# several runs - "w" / "alter@walterkwaninsurance" / "." / "com" - but Find
# & Replace operates on the flattened story text, so it can be targeted as
# one contiguous string). The commit changes the mailbox name from "walter"
# to "info", giving "info@walterkwaninsurance.com".

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "walter@walterkwaninsurance",  # FindText
    $true,                         # MatchCase
    $false,                        # MatchWholeWord
    $false,                        # MatchWildcards
    $false,                        # MatchSoundsLike
    $false,                        # MatchAllWordForms
    $true,                         # Forward
    1,                             # Wrap (wdFindContinue)
    $false,                        # Format
    "info@walterkwaninsurance",    # ReplaceWith
    2                              # Replace (wdReplaceAll)
)
